$wb = $excel.ActiveWorkbook

# --- Sheet1: update "execute" for loginlogout test from "no" to "Yes" ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C2").Value = "Yes"
$ws1.Range("E2").Select()

# --- DATA sheet: update "execute" column to "No" for rows 3-6, and fix username in row 6 ---
$ws2 = $wb.Worksheets.Item("DATA")
$ws2.Range("B3").Value = "No"
$ws2.Range("B4").Value = "No"
$ws2.Range("B5").Value = "No"
$ws2.Range("B6").Value = "No"
$ws2.Range("D6").Value = "Admin"

# browser column (C2:C6) gets re-entered with a quote-prefix (text forced) style
$ws2.Range("C2").Value = "'chrome"
$ws2.Range("C3").Value = "'chrome"
$ws2.Range("C4").Value = "'chrome"
$ws2.Range("C5").Value = "'chrome"
$ws2.Range("C6").Value = "'chrome"

$ws2.Range("C6").Select()
